$d = $word.ActiveDocument

function Escape-Xml([string]$s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

function Add-BoldGreenOK([string]$anchorText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $rng.Collapse(0)
    $rng.InsertAfter(" OK")
    $rng.Font.Bold = $true
    $rng.HighlightColorIndex = 4
}

# 1) "Commande des Servomoteurs ... (Sophie & Adèle)" -> append bold/green " OK"
Add-BoldGreenOK "(Sophie & Adèle)"

# 2) "Commande génération de Son ... (Eddy, Justine et Sophie)" -> append bold/green " OK"
Add-BoldGreenOK "(Eddy, Justine et Sophie)"

# 3) "Mesure du courant consommé ... (Maxime)" -> append bold/green " OK"
Add-BoldGreenOK "(Maxime)"

# 4) "Déplacement évolué (" / bookmark / "& dégradé ?)" / " THIBAUT & CAPUCINE"
#    -> re-split the (unchanged) run text so the run boundaries -- and the
#       _GoBack bookmark sitting between them -- land right after
#       "...dégradé ?) THI" instead of right after the opening "(".
#    The surrounding text itself is not changed, only where it is cut into
#    runs, so pull the live text out of the document instead of retyping it
#    (keeps accents / the non-breaking space between "dégradé" and "?)"
#    byte-for-byte identical).
$rng = $d.Content
$found = $rng.Find.Execute("Déplacement évolué", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph not found: Déplacement évolué"
}
$para = $rng.Paragraphs(1)
$prng = $para.Range
# Exclude the trailing paragraph mark from the replacement range.
$textRng = $d.Range($prng.Start, $prng.End - 1)
$full = $textRng.Text

$cut1 = $full.IndexOf("?)") + 2      # end of "...dégradé ?)"
$cut2 = $full.IndexOf("THI") + 3     # end of "...THI"

$seg1 = $full.Substring(0, $cut1)
$seg2 = $full.Substring($cut1, $cut2 - $cut1)
$seg3 = $full.Substring($cut2)

$seg1x = Escape-Xml $seg1
$seg2x = Escape-Xml $seg2
$seg3x = Escape-Xml $seg3

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
       '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>' + $seg1x + '</w:t></w:r>' +
       '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">' + $seg2x + '</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
       '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>' + $seg3x + '</w:t></w:r>' +
       '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$textRng.InsertXML($xml)
